$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5 and 6 hold two price records for the same market/product that
# differ only by date/volume/price/origin details. The underlying data
# got re-sorted (weekly update), so row 5 and row 6 swap their
# date / volume / price / unit / origin / $-per-kg / kg-per-unit values.

# New row 5 (was row 6's data)
$ws.Range("D5").Value = 44334
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 12500
$ws.Range("Q5").Value = "$/caja 12 kilos empedrada"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1042
$ws.Range("T5").Value = 12

# New row 6 (was row 5's data)
$ws.Range("D6").Value = 44330
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = "$/caja 18 kilos granel"
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value = 861
$ws.Range("T6").Value = 18
